# Update countries & provincias Spain
# Refresh the "Pais" COVID-19 stats sheet: re-sort a handful of countries
# into their correct rank-by-"Casos totales" position and refresh the
# B:H figures (and timestamp in A1) to the 07:12 data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 27 de Marzo de 2020 a las 07:12'

$ws.Range('A22').Value = 'Israel'
$ws.Range('B22').Value = 3035
$ws.Range('C22').Value = 342
$ws.Range('D22').Value = 79
$ws.Range('E22').Value = 2946
$ws.Range('F22').Value = 49
$ws.Range('G22').Value = 2
$ws.Range('H22').Value = 10

$ws.Range('A23').Value = 'Brasil'
$ws.Range('B23').Value = 2985
$ws.Range('D23').Value = 6
$ws.Range('E23').Value = 2902
$ws.Range('F23').Value = 296

$ws.Range('A24').Value = 'Suecia'
$ws.Range('B24').Value = 2840
$ws.Range('D24').Value = 16
$ws.Range('E24').Value = 2747
$ws.Range('F24').Value = 200
$ws.Range('H24').Value = 77

$ws.Range('B34').Value = 1203
$ws.Range('C34').Value = 2
$ws.Range('E34').Value = 1173

$ws.Range('B35').Value = 1136
$ws.Range('C35').Value = 91
$ws.Range('D35').Value = 97
$ws.Range('E35').Value = 1034
$ws.Range('F35').Value = 11
$ws.Range('G35').Value = 1
$ws.Range('H35').Value = 5

$ws.Range('B44').Value = 745
$ws.Range('C44').Value = 18
$ws.Range('E44').Value = 659

$ws.Range('A67').Value = 'Hungria'
$ws.Range('B67').Value = 300
$ws.Range('C67').Value = 39
$ws.Range('D67').Value = 34
$ws.Range('E67').Value = 256
$ws.Range('F67').Value = 6
$ws.Range('H67').Value = 10

$ws.Range('A68').Value = 'Lituania'
$ws.Range('B68').Value = 299
$ws.Range('D68').Value = 1
$ws.Range('E68').Value = 294
$ws.Range('F68').Value = 1
$ws.Range('H68').Value = 4

$ws.Range('A69').Value = 'Armenia'
$ws.Range('B69').Value = 290
$ws.Range('D69').Value = 18
$ws.Range('E69').Value = 271
$ws.Range('F69').Value = 6
$ws.Range('H69').Value = 1

$ws.Range('A70').Value = 'Marruecos'
$ws.Range('B70').Value = 275
$ws.Range('E70').Value = 256
$ws.Range('F70').Value = 1
$ws.Range('H70').Value = 11

$ws.Range('A71').Value = 'Taiwan'
$ws.Range('B71').Value = 267
$ws.Range('C71').Value = 15
$ws.Range('D71').Value = 30
$ws.Range('E71').Value = 235
$ws.Range('F71').Value = 0
$ws.Range('H71').Value = 2

$ws.Range('A72').Value = 'Bulgaria'
$ws.Range('B72').Value = 264
$ws.Range('D72').Value = 8
$ws.Range('E72').Value = 253
$ws.Range('F72').Value = 8
$ws.Range('H72').Value = 3

$ws.Range('A107').Value = 'Uzbekistan'
$ws.Range('B107').Value = 83
$ws.Range('C107').Value = 8
$ws.Range('E107').Value = 83
$ws.Range('F107').Value = 7
$ws.Range('H107').Value = 0

$ws.Range('A108').Value = 'Martinica'
$ws.Range('E108').Value = 80
$ws.Range('F108').Value = 12
$ws.Range('H108').Value = 1

$ws.Range('A109').Value = 'Mauricio'
$ws.Range('B109').Value = 81
$ws.Range('C109').Value = 0
$ws.Range('D109').Value = 0
$ws.Range('H109').Value = 2

$ws.Range('A110').Value = 'Georgia'
$ws.Range('B110').Value = 79
$ws.Range('D110').Value = 11
$ws.Range('E110').Value = 68
$ws.Range('F110').Value = 1

$ws.Range('A122').Value = 'Banglades'
$ws.Range('B122').Value = 48
$ws.Range('C122').Value = 4
$ws.Range('D122').Value = 11
$ws.Range('E122').Value = 32
$ws.Range('F122').Value = 1
$ws.Range('H122').Value = 5

$ws.Range('A123').Value = 'Kirguistan'
$ws.Range('D123').Value = 0
$ws.Range('E123').Value = 44
$ws.Range('F123').Value = 0
$ws.Range('H123').Value = 0

$ws.Range('A134').Value = 'Isla de Man'
$ws.Range('D134').Value = 0
$ws.Range('E134').Value = 26
$ws.Range('H134').Value = 0

$ws.Range('A135').Value = 'Jamaica'
$ws.Range('B135').Value = 26
$ws.Range('D135').Value = 2
$ws.Range('E135').Value = 23
$ws.Range('H135').Value = 1

$ws.Range('A151').Value = 'Republica de Yibuti'

$ws.Range('A152').Value = 'San Martin (Parte Francesa)'

$ws.Range('A160').Value = 'Granada'

$ws.Range('A161').Value = 'Seychelles'

$ws.Range('A162').Value = 'Antigua y Barbuda'

$ws.Range('A163').Value = 'Mozambique'

$ws.Range('A166').Value = 'Suazilandia'

$ws.Range('A168').Value = 'Laos'

$ws.Range('A171').Value = 'Birmania'

$ws.Range('A172').Value = 'Fiyi'

$ws.Range('A173').Value = 'Montserrat'

$ws.Range('A175').Value = 'Guyana'

$ws.Range('A176').Value = 'Cabo Verde'

$ws.Range('A177').Value = 'Congo'

$ws.Range('A178').Value = 'Santa Sede'

$ws.Range('A179').Value = 'Guinea'

$ws.Range('A180').Value = 'Angola'

$ws.Range('A181').Value = 'Mali'

$ws.Range('A182').Value = 'Liberia'

$ws.Range('A183').Value = 'Butan'
$ws.Range('C183').Value = 1

$ws.Range('A184').Value = 'Mauritania'

$ws.Range('A186').Value = 'Republica de Africa Central'

$ws.Range('A187').Value = 'San Martin (Parte Holandesa)'

$ws.Range('A188').Value = 'San Bartolome'
$ws.Range('C188').Value = 0

$ws.Range('A189').Value = 'Gambia'
$ws.Range('D189').Value = 0
$ws.Range('H189').Value = 1

$ws.Range('A190').Value = 'Zimbabue'
$ws.Range('D190').Value = 0
$ws.Range('H190').Value = 1

$ws.Range('A191').Value = 'Sudan'

$ws.Range('A192').Value = 'Nepal'
$ws.Range('D192').Value = 1
$ws.Range('H192').Value = 0

$ws.Range('A193').Value = 'Santa Lucia'
$ws.Range('D193').Value = 1
$ws.Range('H193').Value = 0

$ws.Range('A194').Value = 'Guinea-Bisau'

$ws.Range('A195').Value = 'Somalia'

$ws.Range('A196').Value = 'Islas Turcas y Caicos'

$ws.Range('A197').Value = 'Islas Virgenes Britanicas'

$ws.Range('A199').Value = 'Anguila'

$ws.Range('A200').Value = 'Belice'

$ws.Range('A202').Value = 'Libia'

$ws.Range('A204').Value = 'San Vicente y las Granadinas'

$ws.Range('A205').Value = 'Papua Nueva Guinea'
